$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "Radio IDs" sheet -> "Radio Labels"
# ---------------------------------------------------------------------------
$wsLabels = $wb.Worksheets.Item("Radio IDs")
$wsLabels.Name = "Radio Labels"

# ---------------------------------------------------------------------------
# 2) Rename data sheets "1".."8" -> "0".."7" (shift every label down by one
#    so they become zero-based indices matching the new "Label Index" column)
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 8; $i++) {
    $ws = $wb.Worksheets.Item("$i")
    $ws.Name = "$($i - 1)"
}

# ---------------------------------------------------------------------------
# 3) Update the "Radio Labels" sheet: header "ID" -> "Label Index", and every
#    numeric label in column A decremented by 1 (1-based IDs -> 0-based index)
# ---------------------------------------------------------------------------
$wsLabels.Range("A1").Value = "Label Index"

$wsLabels.Range("A2").Value = 0
$wsLabels.Range("A3").Value = 1
$wsLabels.Range("A4").Value = 2
$wsLabels.Range("A5").Value = 3
$wsLabels.Range("A6").Value = 4
$wsLabels.Range("A7").Value = 5
$wsLabels.Range("A8").Value = 6
$wsLabels.Range("A9").Value = 7
$wsLabels.Range("A12").Value = 8

$wsLabels.Range("A13").Select()

# ---------------------------------------------------------------------------
# 4) Sheet "0" (originally "1") data got sorted ascending by Noise Amplitude
#    (column A), which is the same as sorting descending by SNR (column B).
# ---------------------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("0")
$sortRange = $ws0.Range("A2:B14")
$sortRange.Sort($ws0.Range("A2:A14"), 1)
$ws0.Range("B18").Select()
